$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert a new row at 38 ("Charge Count Quantity") - format is copied from
#    row 37 automatically by Insert().
# ---------------------------------------------------------------------------
$ws.Rows.Item(38).Insert()
$ws.Range("B38").Value = "Charge Count Quantity"
$ws.Range("E38").Value = "/cq-res-doc:CustodyQueryResults/cq-res-ext:Custody/j:Charge[@structures:id=/cq-res-doc:CustodyQueryResults/cq-res-ext:Custody/j:ActivityChargeAssociation/j:Charge/@structures:ref]/j:ChargeCountQuantity"
$ws.Rows.Item(38).RowHeight = 42

# ---------------------------------------------------------------------------
# 2. Insert a new row at 41 ("Statute Section ID" / "Statute or Ordinance
#    Section Number") - format copied from row 40 (formerly old row 39).
# ---------------------------------------------------------------------------
$ws.Rows.Item(41).Insert()
$ws.Range("C41").Value = "Statute Section ID"
$ws.Range("B41").Value = "Statute or Ordinance Section Number"
$ws.Rows.Item(41).RowHeight = 56

# ---------------------------------------------------------------------------
# 3. Row 40 (formerly old row 39) gets a new XPath mapping in column E - the
#    old XPath text that used to live here moves down to the new row 41.
# ---------------------------------------------------------------------------
$ws.Range("E41").Value = $ws.Range("E40").Value2
$ws.Range("E40").Value = "/cq-res-doc:CustodyQueryResults/cq-res-ext:Custody/j:Charge[@structures:id=/cq-res-doc:CustodyQueryResults/cq-res-ext:Custody/j:ActivityChargeAssociation/j:Charge/@structures:ref]/j:ChargeStatute/j:StatuteCodeIdentification/nc:IdentificationID"

# ---------------------------------------------------------------------------
# 4. Remove the two blank placeholder rows that are now at 60/61 (they used
#    to sit right after the "Don't Use" block before the shift).
# ---------------------------------------------------------------------------
$ws.Rows.Item(60).Delete()
$ws.Rows.Item(60).Delete()

# Row 61 (formerly the un-hidden blank row that followed) now needs to be
# hidden.
$ws.Rows.Item(61).Hidden = $true

# ---------------------------------------------------------------------------
# 5. Toggle hidden state on a couple of placeholder row blocks further down.
# ---------------------------------------------------------------------------
$ws.Rows.Item(273).Hidden = $false
$ws.Rows.Item(274).Hidden = $false
$ws.Rows.Item(287).Hidden = $true
$ws.Rows.Item(288).Hidden = $true

# ---------------------------------------------------------------------------
# 6. Append two new blank rows at the bottom (385 / 386).
# ---------------------------------------------------------------------------
$ws.Rows.Item(385).Insert()
$ws.Rows.Item(386).Insert()

# ---------------------------------------------------------------------------
# 7. Update the view: scroll frozen pane back to the top and select C38.
# ---------------------------------------------------------------------------
$ws.Range("C38").Select()
